# Edit: regenerate the "count" values (column F) in the StaffingReal sheet,
# corresponding to an improved query against available_shift, and adjust the
# book-view tab ratio. This mirrors lib/tasks/data/staffing_real.xlsx diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StaffingReal")

# --- Book view cosmetic change: tabRatio 992 -> 985 -------------------------
try {
    $win = $excel.ActiveWindow
    $win.TabRatio = 0.985
} catch {
    # TabRatio may not be settable in every host; ignore if unsupported.
}

# --- Data changes: column F ("count") for the listed rows ------------------
$rows = @(92,93,94,95,96,97,98,99,100,101,102,103,104,105,106,107,108,109,110,111,112,113,115,116,117,118,119,120,121,122,123,124,125,126,127,128,129,130,131,132,133,134,135,136,137,138,139,140,141,142,143,144,145,146,147,148,149,150,151,152,162,206,273,274,275,276,277,278,279,280,281,282,283,284,285,286,287,288,289,290,291,292,293,294,295,296,297,298,299,300,301,302,303)
$vals = @(59,83,111,112,121,60,61,87,100,95,96,99,131,97,105,53,123,82,78,130,94,94,127,48,98,128,96,135,99,51,70,93,95,93,50,54,70,70,70,76,82,104,85,87,43,96,65,69,100,71,54,62,94,41,82,105,81,83,99,82,70,70,70,70,93,95,93,50,54,70,88,21,76,82,104,85,87,43,96,65,69,100,71,54,62,94,41,82,105,81,83,99,82)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 6).Value = $vals[$i]
}
